# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell($cell, $text) {
    # Assign as plain text and restore the default "Normal" style so
    # numeric-looking strings (e.g. "1.00") are not silently coerced
    # into numbers by Excel's usual text-to-number inference, while
    # keeping cell formatting identical to the untouched cells.
    $ws.Range($cell).Value = "'" + $text
    $ws.Range($cell).Style = "Normal"
}

Set-TextCell "D2" "59.398.88"
Set-TextCell "E2" "  +0.63%  "
Set-TextCell "D3" "2.599.68"
Set-TextCell "E3" "  +0.44%  "
Set-TextCell "E4" "  -0.06%  "
Set-TextCell "D5" "536.37"
Set-TextCell "E5" "  +2.52%  "
Set-TextCell "D6" "141.20"
Set-TextCell "E6" "  +1.37%  "
Set-TextCell "E7" "  +0.12%  "
Set-TextCell "E8" "  +0.20%  "
Set-TextCell "D9" "6.50"
Set-TextCell "E9" "  -0.87%  "
Set-TextCell "E10" "  +1.29%  "
Set-TextCell "E11" "  +1.52%  "
Set-TextCell "E12" "  -0.74%  "
Set-TextCell "D13" "3.056.72"
Set-TextCell "E13" "  +0.50%  "
Set-TextCell "D14" "59.299.57"
Set-TextCell "E14" "  +0.55%  "
Set-TextCell "E15" "  +1.10%  "
Set-TextCell "D16" "2.611.51"
Set-TextCell "E16" "  +0.97%  "
Set-TextCell "E17" "  +0.27%  "
Set-TextCell "D18" "340.85"
Set-TextCell "E18" "  +0.67%  "
Set-TextCell "E19" "  +1.35%  "
Set-TextCell "D20" "10.09"
Set-TextCell "E20" "  +0.02%  "
Set-TextCell "D21" "6.36"
Set-TextCell "E21" "  -2.09%  "
Set-TextCell "D23" "67.46"
Set-TextCell "E23" "  +2.05%  "
Set-TextCell "E24" "  +1.23%  "
Set-TextCell "E25" "  -1.60%  "
Set-TextCell "E26" "  +0.09%  "
Set-TextCell "D27" "7.22"
Set-TextCell "E27" "  +2.94%  "
Set-TextCell "E28" "  +2.41%  "
Set-TextCell "E29" "  +0.01%  "
Set-TextCell "E30" "  +4.97%  "
Set-TextCell "E31" "  -1.87%  "
Set-TextCell "D32" "18.81"
Set-TextCell "E32" "  +0.65%  "
Set-TextCell "D33" "150.06"
Set-TextCell "E33" "  +0.61%  "
Set-TextCell "D34" "3.97"
Set-TextCell "E34" "  -0.50%  "
Set-TextCell "E35" "  -0.98%  "
Set-TextCell "E36" "  -0.38%  "
Set-TextCell "D37" "0.831"
Set-TextCell "E37" "  +2.03%  "
Set-TextCell "D38" "0.824"
Set-TextCell "E38" "  -0.35%  "
Set-TextCell "E39" "  +0.51%  "
Set-TextCell "D40" "1.00"
Set-TextCell "E40" "  +0.23%  "
Set-TextCell "D41" "272.29"
Set-TextCell "E41" "  +0.01%  "
Set-TextCell "D42" "0.598"
Set-TextCell "E42" "  +1.56%  "
Set-TextCell "D43" "10.73"
Set-TextCell "E43" "  -0.15%  "
Set-TextCell "D44" "0.0953"
Set-TextCell "E44" "  -0.09%  "
Set-TextCell "E45" "  +1.05%  "
Set-TextCell "D46" "18.59"
Set-TextCell "E46" "  +3.65%  "
Set-TextCell "D47" "1.944.31"
Set-TextCell "E47" "  -1.02%  "
Set-TextCell "E48" "  +1.55%  "
Set-TextCell "E49" "  -0.47%  "
Set-TextCell "D50" "111.42"
Set-TextCell "E50" "  -1.54%  "
Set-TextCell "E51" "  +0.36%  "
